# modeleJury mis a jour
#
# The "JURY" header block in rows 2-5 shows four stacked title lines:
#   F2 = "BUT 2 INFORMATIQUE"
#   F3 = "SEMESTRE 2"
#   F4 = "2023 -2024"      <- placeholder year, to be cleared
#   F5 = "JURY DU <date>"  <- placeholder date, to be cleared
#
# This update clears the two placeholder cells (F4/F5) so the template no
# longer carries a stale/hard-coded year and date, while keeping their
# existing cell formatting (style) intact. It also leaves the active
# selection on F21, matching where the editor's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stale "2023 -2024" year placeholder, keep formatting.
$ws.Range("F4").ClearContents()

# Clear the stale "JURY DU <date>" date placeholder, keep formatting.
$ws.Range("F5").ClearContents()

# Move/save the selection to F21, matching the saved workbook state.
$ws.Range("F21").Select()
